# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "91.237.76"
$ws.Range("E2").Value = "  +4.45%  "
$ws.Range("D3").Value = "3.205.68"
$ws.Range("E3").Value = "  +1.51%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "220.27"
$ws.Range("E5").Value = "  +6.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "637.18"
$ws.Range("E6").Value = "  +5.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.403"
$ws.Range("E7").Value = "  +6.13%  "
$ws.Range("E8").Value = "  +7.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").Value = "3.207.20"
$ws.Range("E10").Value = "  +1.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.575"
$ws.Range("E11").Value = "  +8.30%  "
$ws.Range("E12").Value = "  +4.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000260"
$ws.Range("E13").Value = "  +7.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.43"
$ws.Range("E14").Value = "  +3.96%  "
$ws.Range("D15").Value = "90.905.07"
$ws.Range("E15").Value = "  +4.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "33.28"
$ws.Range("E16").Value = "  +4.22%  "
$ws.Range("D17").Value = "3.794.95"
$ws.Range("E17").Value = "  +1.45%  "
$ws.Range("D18").Value = "3.206.04"
$ws.Range("E18").Value = "  +1.37%  "
$ws.Range("E19").Value = "  +78.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.34"
$ws.Range("E20").Value = "  +5.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "441.43"
$ws.Range("E21").Value = "  +7.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.45"
$ws.Range("E22").Value = "  +1.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.63"
$ws.Range("E23").Value = "  +2.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.06"
$ws.Range("E24").Value = "  +1.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.33"
$ws.Range("E25").Value = "  +4.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.86"
$ws.Range("E26").Value = "  -0.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "81.28"
$ws.Range("E27").Value = "  +11.48%  "
$ws.Range("D28").Value = "3.376.69"
$ws.Range("E28").Value = "  +1.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.161"
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.22"
$ws.Range("E32").Value = "  +40.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.45"
$ws.Range("E33").Value = "  +3.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "534.34"
$ws.Range("E34").Value = "  -1.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.10"
$ws.Range("E35").Value = "  +6.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.91"
$ws.Range("E36").Value = "  +3.85%  "
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "22.57"
$ws.Range("E38").Value = "  +4.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "22.37"
$ws.Range("E39").Value = "  +2.53%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.128"
$ws.Range("E40").Value = "  -3.09%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("E42").Value = "  +2.98%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.374"
$ws.Range("E44").Value = "  +2.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "146.96"
$ws.Range("E45").Value = "  -1.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "44.62"
$ws.Range("E46").Value = "  +3.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "173.63"
$ws.Range("E47").Value = "  +0.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.126"
$ws.Range("E48").Value = "  +2.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.756"
$ws.Range("E49").Value = "  +9.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.11"
$ws.Range("E50").Value = "  +7.33%  "
$ws.Range("E51").Value = "  +1.87%  "

# Reset number format style on cells we forced to text, to avoid leaving a visible style change
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
